$d = $word.ActiveDocument

function Set-FoundText($doc, $searchText, $newText) {
    $r = $doc.Content
    $found = $r.Find.Execute($searchText)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $r.Text = $newText
    return $r
}

# 1. Title
Set-FoundText $d "Untangling the Code of Time" "The Significance of Art: A Reflection of Our Humanity" | Out-Null

# 2. Author name: "Alex Smith" -> three runs: "Dr" / "." / " Ava Roberts"
$r = $d.Content
$r.Find.Execute("Alex Smith") | Out-Null
$r.Text = "Dr"
$insertPoint = $r.End
$r2 = $d.Range($insertPoint, $insertPoint)
$r2.InsertAfter(".")
$insertPoint2 = $r2.End
$r3 = $d.Range($insertPoint2, $insertPoint2)
$r3.InsertAfter(" Ava Roberts")

# 3. Email user name and domain
Set-FoundText $d "alexsmith@quantumtheory" "aroberts@excellentedu" | Out-Null
Set-FoundText $d "net" "org" | Out-Null

# 4. Body paragraph - three chunks separated by <w:br/> line breaks.
Set-FoundText $d "Within the vast tapestry of the universe, time stands as an enigma, an entity that eludes our grasp yet governs our existence" "In the realm of human experience, art stands as a testament to our creative spirit and capacity for expression" | Out-Null
Set-FoundText $d " Across diverse disciplines, scientific minds and inquiring spirits seek to unravel the intricate code that weaves the fabric of time" " Throughout history, art has served as a vessel for emotions, ideas, and stories, weaving together the tapestry of our shared humanity" | Out-Null
Set-FoundText $d " From the theories of theoretical physics to the musings of philosophy and art, our pursuit of understanding time's enigmatic nature reveals a tapestry of theories, conjectures, and intellectual explorations that span the ages" " Whether it's the intricate strokes of a Renaissance painting, the soaring melodies of a classical symphony, or the poignant words of a carefully crafted poem, art has the power to transport us to different times, places, and perspectives, enriching our understanding of the world around us" | Out-Null

Set-FoundText $d "In the realm of theoretical physics, the mysteries of time have drawn scientists into an enigmatic realm" "In its myriad forms, art reflects our deepest longings and aspirations" | Out-Null
Set-FoundText $d " Through the probing lens of general relativity, time emerges as a malleable aspect of spacetime, affected by the gravitational fields it permeates, distorting and stretching with the dance of celestial bodies" " It captures the essence of joy and sorrow, love and loss, triumph and despair, delving into the complexities of the human condition" | Out-Null

# This sentence grows from one run into two runs (new sentence added) in the target.
$r = $d.Content
$r.Find.Execute(" Quantum mechanics, on the other hand, presents a far more perplexing narrative, where time appears to be an emergent framework, intertwined with the fundamental constituents of matter and energy") | Out-Null
$r.Text = " Through art, we find solace and inspiration, connection and community"
$insertPoint = $r.End
$r2 = $d.Range($insertPoint, $insertPoint)
$r2.InsertAfter(".")
$insertPoint2 = $r2.End
$r3 = $d.Range($insertPoint2, $insertPoint2)
$r3.InsertAfter(" It unveils the hidden depths within us, expanding our horizons and challenging our preconceptions, encouraging us to reflect upon our own existence and place in the universe")

Set-FoundText $d "Beyond the boundaries of science, the artistic and philosophical worlds have embraced time as a muse, a source of profound inspiration" "Furthermore, art fosters critical thinking and problem-solving skills, stimulating creativity and innovation, and enhancing our ability to communicate and collaborate effectively" | Out-Null

# The remaining three sentences (literature/music/visual arts) collapse into a single replacement sentence.
$r = $d.Content
$r.Find.Execute(" In literature, time transforms into an ethereal character, an invisible force that shapes narratives, sculpting the destinies of fictional characters") | Out-Null
$startPos = $r.Start
$r2 = $d.Content
$r2.Find.Execute(" Visual arts, with their frozen moments and transient installations, grant us a glimpse into the fragmented nature of time, capturing fleeting moments in a tangible form") | Out-Null
$endPos = $r2.End
$rWhole = $d.Range($startPos, $endPos)
$rWhole.Text = " By engaging with art, we develop our analytical and interpretive abilities, learning to decipher symbols, identify patterns, and appreciate subtleties, skills that are invaluable in any field of endeavor"

# 5. Summary paragraph
Set-FoundText $d "Our quest to decipher the code of time is a journey that reaches beyond the confines of individual disciplines, merging the realms of science, arts, and philosophy" "In essence, art is a fundamental aspect of human existence, reflecting our emotions, ideas, and experiences, enriching our understanding of ourselves, others, and the world around us" | Out-Null
Set-FoundText $d " The theories of general relativity and quantum mechanics provide scientific frameworks for contemplating time's malleability and emergent properties" " Embracing art in its various forms cultivates creativity, critical thinking, communication, and collaboration skills, while inspiring us to reflect upon our own existence" | Out-Null

# The remaining sentences (Meanwhile.../In essence.../humanity's...) collapse into one new sentence.
$r = $d.Content
$r.Find.Execute(" Meanwhile, art and literature explore the subjective, emotional dimensions of time, revealing its power to shape narratives and evoke introspection") | Out-Null
$startPos = $r.Start
$r2 = $d.Content
$r2.Find.Execute("humanity's unrelenting desire to unravel the secrets of the universe and comprehend our place within its intricate tapestry") | Out-Null
$endPos = $r2.End
$rWhole = $d.Range($startPos, $endPos)
$rWhole.Text = " As we delve deeper into the realm of art, we unlock the treasures of our humanity, fostering a more compassionate and interconnected global community"

# 6. Add a new empty paragraph at the very end of the document.
$d.Content.InsertParagraphAfter()

Write-Host "Edit complete"
